# The sheet originally has 5 data rows (rows 2-6). The final state keeps only
# one data row, built from the old row 5's content, with date_of_lab corrected
# from "2024-06-25" to "2024-04-08" (both in column E and inside the JSON
# "context" column). All other former rows (old rows 2, 3, 4, 6) are removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the rows that won't survive in the final sheet (old rows 3,4,5,6 as
# seen from the top; row 2 will be overwritten in place with old row 5's data).
$ws.Range("A3:G6").EntireRow.Delete()

# Rewrite row 2 with the surviving record's data (previously row 5), fixing
# the date_of_lab value along the way.
$ws.Range("A2").Value = "2024-06-25_00:00:00.000_Progress_Notes_91596"
$ws.Range("B2").Value = "<0.06 mg/dL"
$ws.Range("C2").Value = "<1.61 mg/dL"

# D2 stays blank (kappa_lambda_ratio is unknown for this record). Touch the
# font property (a no-op) so the cell is retained in the sheet as an empty
# cell rather than being dropped entirely.
$ws.Range("D2").Value = ""
$ws.Range("D2").Font.Bold = $false

# E2 holds the corrected lab date. Force text formatting first so Excel does
# not reinterpret the string as a date serial number, then restore the
# default "Normal" style so no stray formatting is left on the cell.
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "2024-04-08"
$ws.Range("E2").Style = "Normal"

$ws.Range("F2").Value = "['Labs from 4/8/2024: Kappa <0.06 mg/dL, Lambda <1.61 mg/dL, SPEP with M-spike 0.3 g/dL, IgG kappa']"

$json = @'
{
  "kappa_flc": "<0.06 mg/dL",
  "lambda_flc": "<1.61 mg/dL",
  "kappa_lambda_ratio": null,
  "date_of_lab": "2024-04-08",
  "evidence_sentences": [
    "Labs from 4/8/2024: Kappa <0.06 mg/dL, Lambda <1.61 mg/dL, SPEP with M-spike 0.3 g/dL, IgG kappa"
  ],
  "source_document": "2024-06-25_00:00:00.000_Progress_Notes_91596"
}
'@
$ws.Range("G2").Value = $json

# Writing the multi-line JSON triggers an automatic row-height bump; re-run
# autofit so the row returns to the sheet's normal (non-custom) height, just
# like the original file.
$ws.Rows.Item(2).AutoFit()
